$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.526.12'
$ws.Range('E2').Value = '  -4.95%  '
$ws.Range('D3').Value = '3.455.97'
$ws.Range('E3').Value = '  -6.55%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.54'
$ws.Range('E5').Value = '  -7.50%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.58'
$ws.Range('E6').Value = '  -8.87%  '
$ws.Range('D7').Value = '3.457.81'
$ws.Range('E7').Value = '  -6.46%  '
$ws.Range('E8').Value = '  +0.12%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.476'
$ws.Range('E9').Value = '  -5.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.137'
$ws.Range('E10').Value = '  -6.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.85'
$ws.Range('E11').Value = '  -4.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.418'
$ws.Range('E12').Value = '  -6.16%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000215'
$ws.Range('E13').Value = '  -7.78%  '
$ws.Range('D14').Value = '4.048.50'
$ws.Range('E14').Value = '  -6.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.14'
$ws.Range('E15').Value = '  -5.27%  '
$ws.Range('D16').Value = '3.450.64'
$ws.Range('E16').Value = '  -7.05%  '
$ws.Range('D17').Value = '66.558.62'
$ws.Range('E17').Value = '  -4.85%  '
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.32'
$ws.Range('E19').Value = '  -3.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.85'
$ws.Range('E20').Value = '  -7.66%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '439.87'
$ws.Range('E21').Value = '  -6.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.93'
$ws.Range('E22').Value = '  -14.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.619'
$ws.Range('E23').Value = '  -5.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '76.90'
$ws.Range('E24').Value = '  -4.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.07%  '
$ws.Range('D26').Value = '3.601.18'
$ws.Range('E26').Value = '  -6.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000123'
$ws.Range('E27').Value = '  -4.58%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.96'
$ws.Range('E28').Value = '  -9.39%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.16'
$ws.Range('E29').Value = '  -11.56%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.49'
$ws.Range('E30').Value = '  -6.38%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.55'
$ws.Range('E31').Value = '  -10.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('E33').Value = '  -4.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.36'
$ws.Range('E34').Value = '  -5.55%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.08'
$ws.Range('E35').Value = '  -7.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.84'
$ws.Range('E36').Value = '  -9.11%  '
$ws.Range('D37').Value = '3.448.79'
$ws.Range('E37').Value = '  -6.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.84'
$ws.Range('E38').Value = '  -7.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '173.68'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.12'
$ws.Range('E42').Value = '  -5.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0857'
$ws.Range('E43').Value = '  -5.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.40'
$ws.Range('E44').Value = '  -8.71%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.875'
$ws.Range('E45').Value = '  -6.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.22'
$ws.Range('E46').Value = '  -3.75%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.46'
$ws.Range('E49').Value = '  -14.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.50'
$ws.Range('E50').Value = '  -4.84%  '
$ws.Range('E51').Value = '  -5.38%  '

# Rows 47 and 48 swapped coins (InjectiveProtocol now ranks above ONDO)
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '26.19'
$ws.Range('E47').Value = '  -10.82%  '

$ws.Range('B48').Value = 'ONDO'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.21'
$ws.Range('E48').Value = '  -4.99%  '
